# update repl links to use replit.com
$p = $ppt.ActivePresentation

# Slide 3 ("span example"): update the repl.it link text to the new replit.com link.
$slide3 = $p.Slides.Item(3)
$spanLinkShape = $slide3.Shapes.Item(2)
$spanLinkShape.TextFrame.TextRange.Text = "https://replit.com/@HylandOutreach/SpanExample"

# Slide 5 ("Div example"): resize/reposition the link rectangle and update its link text.
$slide5 = $p.Slides.Item(5)
$divLinkShape = $slide5.Shapes.Item(3)
$divLinkShape.Left = 381
$divLinkShape.Top = 108
$divLinkShape.Width = 522
$divLinkShape.Height = 261
$divLinkShape.TextFrame.TextRange.Text = "https://replit.com/@HylandOutreach/DivExample"
